$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSdnn = @{
    2  = "29.871379591073207"
    3  = "23.734258537001377"
    4  = "19.9648155496385"
    5  = "18.628600778943706"
    6  = "16.2790493823268"
    7  = "19.046189094326152"
    8  = "18.466857897712615"
    9  = "21.45385155342088"
    10 = "22.204304741127483"
    11 = "18.06993412458282"
    12 = "15.48470001172583"
    13 = "16.336820679251357"
    14 = "21.267563599808607"
    15 = "19.996736627512522"
    16 = "10.676169148366048"
    17 = "24.44593306808824"
    18 = "12.300549447255673"
    19 = "19.109729235928555"
    20 = "15.36459887606487"
    21 = "19.160646719617016"
    22 = "14.481527814310157"
    23 = "17.20350490714225"
    24 = "16.24465216380139"
    25 = "21.623741553475075"
    26 = "21.43284376748338"
    27 = "22.421531578161254"
    28 = "14.521597979402125"
}

foreach ($row in $newSdnn.Keys) {
    $text = "ReturnTuple(sdnn=" + $newSdnn[$row] + ")"
    foreach ($col in @("C", "D", "E", "F")) {
        $ws.Range("$col$row").Value = $text
    }
}
